# Update crypto price/volume figures to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.125.59"
$ws.Range("E2").Value = "  +1.16%  "
# Row 3
$ws.Range("D3").Value = "2.648.47"
$ws.Range("E3").Value = "  -0.44%  "
# Row 4
$ws.Range("E4").Value = "  +0.34%  "
# Row 5
$ws.Range("D5").Value = "'607.55"
$ws.Range("E5").Value = "  -0.40%  "
# Row 6
$ws.Range("D6").Value = "'148.26"
$ws.Range("E6").Value = "  +3.38%  "
# Row 7
$ws.Range("E7").Value = "  +0.32%  "
# Row 8
$ws.Range("E8").Value = "  +0.69%  "
# Row 9
$ws.Range("E9").Value = "  +1.81%  "
# Row 10
$ws.Range("E10").Value = "  +6.84%  "
# Row 11
$ws.Range("E11").Value = "  -0.03%  "
# Row 12
$ws.Range("E12").Value = "  -0.89%  "
# Row 13
$ws.Range("D13").Value = "'27.58"
$ws.Range("E13").Value = "  +1.04%  "
# Row 14
$ws.Range("D14").Value = "3.127.96"
$ws.Range("E14").Value = "  +0.12%  "
# Row 15
$ws.Range("D15").Value = "64.036.13"
$ws.Range("E15").Value = "  +1.58%  "
# Row 16
$ws.Range("E16").Value = "  +2.21%  "
# Row 17
$ws.Range("D17").Value = "2.653.15"
$ws.Range("E17").Value = "  +0.00%  "
# Row 18
$ws.Range("D18").Value = "'11.96"
$ws.Range("E18").Value = "  +4.43%  "
# Row 19
$ws.Range("D19").Value = "'4.58"
$ws.Range("E19").Value = "  +3.76%  "
# Row 20
$ws.Range("D20").Value = "'346.20"
$ws.Range("E20").Value = "  +1.24%  "
# Row 21
$ws.Range("D21").Value = "'6.91"
$ws.Range("E21").Value = "  +0.72%  "
# Row 22
$ws.Range("E22").Value = "  +0.17%  "
# Row 23
$ws.Range("E23").Value = "  -0.83%  "
# Row 24
$ws.Range("D24").Value = "'66.24"
$ws.Range("E24").Value = "  -0.88%  "
# Row 25
$ws.Range("E25").Value = "  +8.93%  "
# Row 26
$ws.Range("E26").Value = "  +4.79%  "
# Row 27
$ws.Range("D27").Value = "'9.32"
$ws.Range("E27").Value = "  +7.72%  "
# Row 28
$ws.Range("D28").Value = "'556.08"
$ws.Range("E28").Value = "  +1.75%  "
# Row 29
$ws.Range("D29").Value = "'8.16"
$ws.Range("E29").Value = "  +4.40%  "
# Row 30
$ws.Range("E30").Value = "  -1.27%  "
# Row 31
$ws.Range("E31").Value = "  +0.13%  "
# Row 32
$ws.Range("E32").Value = "  +1.07%  "
# Row 33
$ws.Range("D33").Value = "0.0₃0849"
$ws.Range("E33").Value = "  +5.16%  "
# Row 34
$ws.Range("E34").Value = "  -1.40%  "
# Row 35
$ws.Range("D35").Value = "'5.34"
$ws.Range("E35").Value = "  +3.67%  "
# Row 36
$ws.Range("D36").Value = "'168.79"
$ws.Range("E36").Value = "  -2.44%  "
# Row 37
$ws.Range("D37").Value = "'0.407"
$ws.Range("E37").Value = "  +0.20%  "
# Row 38
$ws.Range("E38").Value = "  -0.03%  "
# Row 39
$ws.Range("E39").Value = "  +5.44%  "
# Row 40
$ws.Range("D40").Value = "'19.36"
$ws.Range("E40").Value = "  +1.10%  "
# Row 41
$ws.Range("E41").Value = "  +0.23%  "
# Row 42
$ws.Range("D42").Value = "'167.05"
$ws.Range("E42").Value = "  -4.70%  "
# Row 43
$ws.Range("D43").Value = "'40.38"
$ws.Range("E43").Value = "  +0.69%  "
# Row 44
$ws.Range("D44").Value = "'3.84"
$ws.Range("E44").Value = "  +2.48%  "
# Row 45
$ws.Range("D45").Value = "'22.08"
$ws.Range("E45").Value = "  -0.62%  "
# Row 46
$ws.Range("E46").Value = "  -0.43%  "
# Row 47
$ws.Range("D47").Value = "'0.628"
$ws.Range("E47").Value = "  -0.52%  "
# Row 48
$ws.Range("D48").Value = "'0.0245"
$ws.Range("E48").Value = "  +2.22%  "
# Row 49
$ws.Range("D49").Value = "'1.98"
$ws.Range("E49").Value = "  +13.68%  "
# Row 50
$ws.Range("D50").Value = "'0.0961"
$ws.Range("E50").Value = "  -0.01%  "
# Row 51
$ws.Range("D51").Value = "'19.03"
$ws.Range("E51").Value = "  +2.22%  "
